$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.045634197174365
$ws.Range("C2").Value = 0.275542105873825
$ws.Range("D2").Value = 0.03046775130748358
$ws.Range("E2").Value = 0.1124537920805614
$ws.Range("F2").Value = 0.7163731155385165
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.6502193167962069
$ws.Range("L2").Value = 0.203420809748259
$ws.Range("M2").Value = 0.2252972918434466
$ws.Range("O2").Value = 2.464290182500179
$ws.Range("B3").Value = 0.9366361524033664
$ws.Range("C3").Value = 0.2575815738641154
$ws.Range("D3").Value = 0.02857376532173816
$ws.Range("E3").Value = 0.1136927504556584
$ws.Range("F3").Value = 0.7173328932413057
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.6606873729560832
$ws.Range("L3").Value = 0.2008182403035761
$ws.Range("M3").Value = 0.2083546590958676
$ws.Range("O3").Value = 2.482262635119895
$ws.Range("B4").Value = 0.8696382720136171
$ws.Range("C4").Value = 0.2464920295403488
$ws.Range("D4").Value = 0.02740426704868071
$ws.Range("E4").Value = 0.1144972361490826
$ws.Range("F4").Value = 0.7184973022593937
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.6675917603935702
$ws.Range("L4").Value = 0.1993230431600566
$ws.Range("M4").Value = 0.197983744319508
$ws.Range("O4").Value = 2.495242534327033
$ws.Range("B5").Value = 0.842319561810541
$ws.Range("C5").Value = 0.2419577229966023
$ws.Range("D5").Value = 0.02692606250322882
$ws.Range("E5").Value = 0.114836088415067
$ws.Range("F5").Value = 0.7191162298727676
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.6705251319483203
$ws.Range("L5").Value = 0.1987396368328831
$ws.Range("M5").Value = 0.1937658041536565
$ws.Range("O5").Value = 2.501020252061338
$ws.Range("B6").Value = 0.8377823638057293
$ws.Range("C6").Value = 0.241203893721746
$ws.Range("D6").Value = 0.02684655980097972
$ws.Range("E6").Value = 0.1148930204349081
$ws.Range("F6").Value = 0.7192277193780043
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.6710194457080725
$ws.Range("L6").Value = 0.1986443283785704
$ws.Range("M6").Value = 0.193065926979969
$ws.Range("O6").Value = 2.502009108945487
$ws.Range("B7").Value = 0.8692699069517857
$ws.Range("C7").Value = 0.2464309395641067
$ws.Range("D7").Value = 0.02739782435103422
$ws.Range("E7").Value = 0.1145017613999479
$ws.Range("F7").Value = 0.7185050648335505
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.6676308361035446
$ws.Range("L7").Value = 0.199315070207625
$ws.Range("M7").Value = 0.197926825780435
$ws.Range("O7").Value = 2.495318478485103
$ws.Range("B8").Value = 1.008067751678709
$ws.Range("C8").Value = 0.2693622887155414
$ws.Range("D8").Value = 0.02981608920282497
$ws.Range("E8").Value = 0.1128719109716995
$ws.Range("F8").Value = 0.7165845669065902
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.6537296161164399
$ws.Range("L8").Value = 0.2025021387521306
$ws.Range("M8").Value = 0.2194489965803541
$ws.Range("O8").Value = 2.470083019476675
$ws.Range("B9").Value = 1.279610971120064
$ws.Range("C9").Value = 0.3138304408676902
$ws.Range("D9").Value = 0.03450493096165985
$ws.Range("E9").Value = 0.1100223379822383
$ws.Range("F9").Value = 0.7173913841942863
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.6302606422624244
$ws.Range("L9").Value = 0.2095659013156919
$ws.Range("M9").Value = 0.26189794110973
$ws.Range("O9").Value = 2.436061116283071
$ws.Range("B10").Value = 1.478660359084074
$ws.Range("C10").Value = 0.3461852634419245
$ws.Range("D10").Value = 0.03791609509146809
$ws.Range("E10").Value = 0.1081390733635839
$ws.Range("F10").Value = 0.7207869119655825
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.6153382926800646
$ws.Range("L10").Value = 0.215250441932298
$ws.Range("M10").Value = 0.2932246992110734
$ws.Range("O10").Value = 2.420542698247374
$ws.Range("B11").Value = 1.569102509158199
$ws.Range("C11").Value = 0.3608336005262629
$ws.Range("D11").Value = 0.03946036856841317
$ws.Range("E11").Value = 0.1073277917433864
$ws.Range("F11").Value = 0.7229435143738812
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.6090554886343504
$ws.Range("L11").Value = 0.2179436861987938
$ws.Range("M11").Value = 0.3075046140601856
$ws.Range("O11").Value = 2.415551732666387
$ws.Range("B12").Value = 1.603333855341759
$ws.Range("C12").Value = 0.3663702309458756
$ws.Range("D12").Value = 0.04004404255876182
$ws.Range("E12").Value = 0.1070270982520678
$ws.Range("F12").Value = 0.7238483971280374
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.6067492080171561
$ws.Range("L12").Value = 0.218978942251411
$ws.Range("M12").Value = 0.3129160324502962
$ws.Range("O12").Value = 2.413960055005731
$ws.Range("B13").Value = 1.595962313461484
$ws.Range("C13").Value = 0.3651782848844221
$ws.Range("D13").Value = 0.03991838775970535
$ws.Range("E13").Value = 0.107091568111964
$ws.Range("F13").Value = 0.7236495870778015
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.6072426622645715
$ws.Range("L13").Value = 0.2187552979632841
$ws.Range("M13").Value = 0.311750416515288
$ws.Range("O13").Value = 2.414289571086073
$ws.Range("B14").Value = 1.571919096420459
$ws.Range("C14").Value = 0.3612893121558329
$ws.Range("D14").Value = 0.03950841018991014
$ws.Range("E14").Value = 0.1073029228746969
$ws.Range("F14").Value = 0.7230161901957359
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.6088642877647352
$ws.Range("L14").Value = 0.2180285493230514
$ws.Range("M14").Value = 0.3079497376315743
$ws.Range("O14").Value = 2.415414801287852
$ws.Range("B15").Value = 1.557189646415168
$ws.Range("C15").Value = 0.3589058436480457
$ws.Range("D15").Value = 0.039257141710479
$ws.Range("E15").Value = 0.1074332327105125
$ws.Range("F15").Value = 0.7226397123606318
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.6098670773259229
$ws.Range("L15").Value = 0.2175853966415815
$ws.Range("M15").Value = 0.3056222149895049
$ws.Range("O15").Value = 2.416142909029617
$ws.Range("B16").Value = 1.472747453689919
$ws.Range("C16").Value = 0.3452265246706929
$ws.Range("D16").Value = 0.03781501980650859
$ws.Range("E16").Value = 0.1081930056654397
$ws.Range("F16").Value = 0.7206583064176542
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.6157590765512637
$ws.Range("L16").Value = 0.2150765873835923
$ws.Range("M16").Value = 0.2922920373932669
$ws.Range("O16").Value = 2.420910564972701
$ws.Range("B17").Value = 1.420916340722727
$ws.Range("C17").Value = 0.3368165494904929
$ws.Range("D17").Value = 0.03692838484843719
$ws.Range("E17").Value = 0.1086707284572481
$ws.Range("F17").Value = 0.7195996765182997
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.6195032256755688
$ws.Range("L17").Value = 0.2135649658058298
$ws.Range("M17").Value = 0.2841216905598642
$ws.Range("O17").Value = 2.42436575363763
$ws.Range("B18").Value = 1.391094547648777
$ws.Range("C18").Value = 0.3319727803984449
$ws.Range("D18").Value = 0.03641771356475942
$ws.Range("E18").Value = 0.1089497782171569
$ws.Range("F18").Value = 0.7190483678258275
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.6217043245301745
$ws.Range("L18").Value = 0.2127056252984971
$ws.Range("M18").Value = 0.2794250894664714
$ws.Range("O18").Value = 2.426547727488753
$ws.Range("B19").Value = 1.380995755283379
$ws.Range("C19").Value = 0.3303316437267654
$ws.Range("D19").Value = 0.03624468942762604
$ws.Range("E19").Value = 0.1090449944511576
$ws.Range("F19").Value = 0.7188715879432763
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.6224577417810018
$ws.Range("L19").Value = 0.2124164042862731
$ws.Range("M19").Value = 0.2778353846108956
$ws.Range("O19").Value = 2.427319909992633
$ws.Range("B20").Value = 1.426434889908364
$ws.Range("C20").Value = 0.337712488296404
$ws.Range("D20").Value = 0.03702284157881053
$ws.Range("E20").Value = 0.1086194315536666
$ws.Range("F20").Value = 0.7197064077779274
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.6190997305470063
$ws.Range("L20").Value = 0.2137248351101988
$ws.Range("M20").Value = 0.2849911536356657
$ws.Range("O20").Value = 2.423977792118563
$ws.Range("B21").Value = 1.578981651589061
$ws.Range("C21").Value = 0.3624318817684582
$ws.Range("D21").Value = 0.03962886089566098
$ws.Range("E21").Value = 0.1072406659582832
$ws.Range("F21").Value = 0.7231998380498652
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.6083859975026726
$ws.Range("L21").Value = 0.2182415959207162
$ws.Range("M21").Value = 0.3090659848330688
$ws.Range("O21").Value = 2.415076191547655
$ws.Range("B22").Value = 1.678579074148104
$ws.Range("C22").Value = 0.3785268095889762
$ws.Range("D22").Value = 0.04132556677728161
$ws.Range("E22").Value = 0.1063775688459541
$ws.Range("F22").Value = 0.7259973379648841
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.6018088509698813
$ws.Range("L22").Value = 0.2212832007074184
$ws.Range("M22").Value = 0.3248230270383488
$ws.Range("O22").Value = 2.41099755145143
$ws.Range("B23").Value = 1.625431870907107
$ws.Range("C23").Value = 0.3699422913206263
$ws.Range("D23").Value = 0.04042060650521506
$ws.Range("E23").Value = 0.1068347460142665
$ws.Range("F23").Value = 0.7244571235114847
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.605280254388834
$ws.Range("L23").Value = 0.2196516531520842
$ws.Range("M23").Value = 0.3164112064830604
$ws.Range("O23").Value = 2.413014992727199
$ws.Range("B24").Value = 1.423940026560388
$ws.Range("C24").Value = 0.3373074616341967
$ws.Range("D24").Value = 0.03698014059519039
$ws.Range("E24").Value = 0.1086426091756401
$ws.Range("F24").Value = 0.7196579760698043
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.6192819995163603
$ws.Range("L24").Value = 0.2136525279460102
$ws.Range("M24").Value = 0.2845980673647688
$ws.Range("O24").Value = 2.424152580487032
$ws.Range("B25").Value = 1.206226402677032
$ws.Range("C25").Value = 0.3018553224590903
$ws.Range("D25").Value = 0.03324231280993217
$ws.Range("E25").Value = 0.1107562205916771
$ws.Range("F25").Value = 0.7166820004923053
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.6362028030439042
$ws.Range("L25").Value = 0.2075679747205825
$ws.Range("M25").Value = 0.2503892155922358
$ws.Range("O25").Value = 2.44360417470611
